$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("K-SVMeans")

# Add the new k = 3 column (C) values for rows 13-17
$ws.Range("C13").Value = 58.23
$ws.Range("C14").Value = 55.67
$ws.Range("C15").Value = 50.26
$ws.Range("C16").Value = 54.89
$ws.Range("C17").Value = 53.24

# Add the AVERAGE formula in C18, matching B18's style
$ws.Range("C18").Formula = "=AVERAGE(C13:C17)"
$ws.Range("B18").Copy()
$ws.Range("C18").PasteSpecial(-4122)  # xlPasteFormats

# Update the selection to reflect the new active cell
$ws.Range("D17").Select()
